# Auto-generated edit applying the Mandragora_Profits market-data refresh
# (scheduled runner price/profit recompute) per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3729.3333
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3729.3333
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3729.3333
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -4079.3333
$ws.Range("H129").Value = 1268
$ws.Range("I129").Value = 587.0909
$ws.Range("J129").Value = 1470.4324
$ws.Range("K129").Value = 1761.2727
$ws.Range("L129").Value = 4411.2972
$ws.Range("M129").Value = 3238.7273
$ws.Range("N129").Value = -14411.2972
$ws.Range("H137").Value = 1966.24
$ws.Range("I137").Value = 2654
$ws.Range("K137").Value = 7962
$ws.Range("M137").Value = -5412

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7519.02
$ws.Range("I32").Value = 6458.911
$ws.Range("K32").Value = 6458.911
$ws.Range("M32").Value = -6171.911
$ws.Range("H74").Value = 2223.6956
$ws.Range("I74").Value = 1841.3889
$ws.Range("J74").Value = 3600
$ws.Range("K74").Value = 1841.3889
$ws.Range("L74").Value = 3600
$ws.Range("M74").Value = -967.3888999999999
$ws.Range("N74").Value = -5348
$ws.Range("H77").Value = 2223.6956
$ws.Range("I77").Value = 1841.3889
$ws.Range("J77").Value = 3600
$ws.Range("K77").Value = 9206.9445
$ws.Range("L77").Value = 18000
$ws.Range("M77").Value = -4838.9445
$ws.Range("N77").Value = -26736
$ws.Range("H122").Value = 1111.8636
$ws.Range("I122").Value = 1025.4
$ws.Range("J122").Value = 1297.1428
$ws.Range("K122").Value = 3076.2
$ws.Range("L122").Value = 3891.4284
$ws.Range("M122").Value = -626.2000000000003
$ws.Range("N122").Value = -8791.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5749368
$ws.Range("I31").Value = 1880.3235
$ws.Range("J31").Value = 13891642
$ws.Range("K31").Value = 1880.3235
$ws.Range("L31").Value = 13891642
$ws.Range("M31").Value = -1585.3235
$ws.Range("N31").Value = -13892232
$ws.Range("H32").Value = 9800
$ws.Range("I32").Value = 8000
$ws.Range("K32").Value = 8000
$ws.Range("M32").Value = -7684
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H34").Value = 5749368
$ws.Range("I34").Value = 1880.3235
$ws.Range("J34").Value = 13891642
$ws.Range("K34").Value = 1880.3235
$ws.Range("L34").Value = 13891642
$ws.Range("M34").Value = -1678.3235
$ws.Range("N34").Value = -13892046
$ws.Range("H36").Value = 9247.5
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 9997.272000000001
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 9997.272000000001
$ws.Range("M36").Value = -612
$ws.Range("N36").Value = -10773.272
$ws.Range("H38").Value = 99454.55
$ws.Range("I38").Value = 7000
$ws.Range("J38").Value = 120000
$ws.Range("K38").Value = 7000
$ws.Range("L38").Value = 120000
$ws.Range("M38").Value = -6623
$ws.Range("N38").Value = -120754
$ws.Range("H39").Value = 10000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 10000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -10782
$ws.Range("H40").Value = 9247.5
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 9997.272000000001
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 9997.272000000001
$ws.Range("M40").Value = -840
$ws.Range("N40").Value = -10317.272
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H45").Value = 9900
$ws.Range("I45").Value = 9000
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 9000
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -8407
$ws.Range("N45").Value = -11186
$ws.Range("H46").Value = 99454.55
$ws.Range("I46").Value = 7000
$ws.Range("J46").Value = 120000
$ws.Range("K46").Value = 7000
$ws.Range("L46").Value = 120000
$ws.Range("M46").Value = -6789
$ws.Range("N46").Value = -120422
$ws.Range("H49").Value = 10000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 10000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 10000
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -10364
$ws.Range("H50").Value = 5500
$ws.Range("I50").Value = 5500
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 5500
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -4875
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 13333.333
$ws.Range("J51").Value = 10000
$ws.Range("L51").Value = 10000
$ws.Range("N51").Value = -11472
$ws.Range("H56").Value = 5000
$ws.Range("I56").Value = 5000
$ws.Range("K56").Value = 5000
$ws.Range("M56").Value = -4155
$ws.Range("H57").Value = 9250
$ws.Range("J57").Value = 9250
$ws.Range("L57").Value = 9250
$ws.Range("N57").Value = -10370
$ws.Range("H61").Value = 13333.333
$ws.Range("J61").Value = 10000
$ws.Range("L61").Value = 10000
$ws.Range("N61").Value = -10696

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 97.25
$ws.Range("I2").Value = 60.5
$ws.Range("J2").Value = 115.625
$ws.Range("K2").Value = 363
$ws.Range("L2").Value = 693.75
$ws.Range("M2").Value = -250
$ws.Range("N2").Value = -919.75
$ws.Range("H136").Value = 4409.778
$ws.Range("I136").Value = 3448.3333
$ws.Range("J136").Value = 6332.6665
$ws.Range("K136").Value = 10344.9999
$ws.Range("L136").Value = 18997.9995
$ws.Range("M136").Value = -5244.999899999999
$ws.Range("N136").Value = -29197.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1667240.4
$ws.Range("I46").Value = 375
$ws.Range("J46").Value = 2500673
$ws.Range("K46").Value = 375
$ws.Range("L46").Value = 2500673
$ws.Range("M46").Value = -187
$ws.Range("N46").Value = -2501049
